# "Generate Report for Handback" — refresh the handback status rows for the
# 3db0f42a-77e3-4615-8e75-beae3b4e46b9 file in both locale report sheets:
#   column E = Correspond Handoff Datetime
#   column H = Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-13 02:45:10"
$zhcn.Range("H3").Value = "2016-03-13 02:45:29"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-13 02:45:14"
$dede.Range("H3").Value = "2016-03-13 02:45:35"
